$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 30,3
$arr[0,0]=-0.4692756063703969; $arr[0,1]=-2.234323765550341; $arr[0,2]=0.2831644466944987
$arr[1,0]=0.5019868420703071; $arr[1,1]=-2.09162780216762; $arr[1,2]=-3.543231316975185
$arr[2,0]=2.689777469634955; $arr[2,1]=-1.050217630181979; $arr[2,2]=-0.7848469325474046
$arr[3,0]=-6.342504692077637; $arr[3,1]=1.490384888648992; $arr[3,2]=3.834645700454709
$arr[4,0]=-5.673670666558469; $arr[4,1]=3.692834442002409; $arr[4,2]=2.466138475707613
$arr[5,0]=3.314056737082347; $arr[5,1]=6.166058949061821; $arr[5,2]=0.5118652531078883
$arr[6,0]=2.383262787546406; $arr[6,1]=-4.818590433469862; $arr[6,2]=0.4682122468948429
$arr[7,0]=8.378727106537124; $arr[7,1]=0.4666124773877008; $arr[7,2]=-1.115668596540174
$arr[8,0]=-0.3396976134606717; $arr[8,1]=0.4775580519012043; $arr[8,2]=1.77553138222013
$arr[9,0]=-5.64267150844849; $arr[9,1]=3.583157830578955; $arr[9,2]=3.719863629341134
$arr[10,0]=-1.766980731487296; $arr[10,1]=16.75916714668267; $arr[10,2]=-0.6191982626914752
$arr[11,0]=1.7845652954919; $arr[11,1]=-3.40247355188642; $arr[11,2]=4.17544686794281
$arr[12,0]=8.421752619743332; $arr[12,1]=-2.540225854941783; $arr[12,2]=-3.182726718698209
$arr[13,0]=0.1595347711017858; $arr[13,1]=-0.04787222487585119; $arr[13,2]=1.229899188450423
$arr[14,0]=-3.870756162064419; $arr[14,1]=2.866504958271982; $arr[14,2]=3.126545447962625
$arr[15,0]=-3.578196597099306; $arr[15,1]=14.92317203113009; $arr[15,2]=-0.7331702096121586
$arr[16,0]=0.1684449655669233; $arr[16,1]=-0.0544045226915344; $arr[16,2]=0.7549586423805885
$arr[17,0]=9.868156909942579; $arr[17,1]=-0.5477856159210392; $arr[17,2]=-1.925824952125532
$arr[18,0]=0.4101533753531434; $arr[18,1]=-0.5063240064041949; $arr[18,2]=1.204618217263898
$arr[19,0]=-4.432378809792649; $arr[19,1]=1.422139526265; $arr[19,2]=3.637724219049721
$arr[20,0]=0.1461163333484263; $arr[20,1]=9.942505402224434; $arr[20,2]=-0.5480119875499323
$arr[21,0]=1.381819248199463; $arr[21,1]=-3.363096782139367; $arr[21,2]=0.9756205422537669
$arr[22,0]=5.626129477364711; $arr[22,1]=-2.49930441464691; $arr[22,2]=-2.115894814900024
$arr[23,0]=1.154986371312843; $arr[23,1]=-0.3653818828718969; $arr[23,2]=1.201646787779655
$arr[24,0]=-3.968483030796053; $arr[24,1]=2.245044922828678; $arr[24,2]=4.008097422122955
$arr[25,0]=-0.9093820844377891; $arr[25,1]=8.189985391071852; $arr[25,2]=-0.7311917849949312
$arr[26,0]=0.1449330214943292; $arr[26,1]=1.165589673178594; $arr[26,2]=-0.2754184092794358
$arr[27,0]=3.144946085555223; $arr[27,1]=-6.092144768578668; $arr[27,2]=-2.065853960173481
$arr[28,0]=2.411846300108097; $arr[28,1]=-2.998377110842354; $arr[28,2]=-1.209750396864759
$arr[29,0]=-1.357690974644261; $arr[29,1]=0.8570090349231524; $arr[29,2]=1.883484615598408

$ws.Range("A2:C31").Value2 = $arr
